# Exercicios Lacos De Repeticao - finish Exercicio 5 / update selection
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Check" box for exercise 5 (row 6, column B) as completed (TRUE)
$ws.Range("B6").Value = $true

# Move/restore the active cell selection to A11 (matches the saved cursor position)
$ws.Range("A11").Select()
